# Issue 43 (csarven): fix "qb:dataset" -> "qb:dataSet" text box on slide 1
# and widen it so the (now slightly longer) label still fits, matching the
# canonical OOXML diff:
#   <a:ext cx="841897" .../>  -> cx="865943"
#   <a:t>qb:dataset</a:t>      -> <a:t>qb:dataSet</a:t>

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 30" -and $candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -eq "qb:dataset") {
            $shp = $candidate
            break
        }
    }
}

if ($shp -eq $null) {
    throw "Could not locate the 'qb:dataset' TextBox 30 shape on slide 1"
}

# Update the run text (single run is preserved since the whole range is
# replaced in one assignment).
$shp.TextFrame.TextRange.Text = "qb:dataSet"

# Widen the textbox (cx 841897 -> 865943 EMU); height (cy) is unchanged.
# Width/Height on the Shape object are expressed in points (1 pt = 12700 EMU:
# 865943 EMU = 68.18448818897637 pt); using the EMU/12700 quotient directly
# rounds (through the engine's internal float32 storage) down to 865942 EMU,
# so a literal point value that lands in the correct float32 bucket is used
# instead to reproduce the exact target EMU value.
$shp.Width = 68.1845
